$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2: ECs -> FAPs (was ECs -> ECs) ---
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 1.738079666666666
$ws.Range("H2").Value = 5.214238999999999
$ws.Range("I2").Value = 0.2081154188575857
$ws.Range("J2").Value = 0.2081154188575857
$ws.Range("M2").Value = 0.1546876666666667
$ws.Range("N2").Value = 0.464063
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.2688594881174444
$ws.Range("R2").Value = 2.419735393057
$ws.Range("S2").Value = 0.2081154188575857
$ws.Range("T2").Value = 0.2081154188575857

# --- Update row 3: FAPs -> FAPs (was ECs -> FAPs) ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 4.080633
$ws.Range("H3").Value = 12.241899
$ws.Range("I3").Value = 0.4886097353798435
$ws.Range("J3").Value = 0.4886097353798435
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.6312235972929999
$ws.Range("R3").Value = 5.681012375637
$ws.Range("S3").Value = 0.4886097353798435
$ws.Range("T3").Value = 0.4886097353798435

# --- Update row 4: MuSCs -> FAPs (was FAPs -> ECs) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 2.532805333333334
$ws.Range("H4").Value = 7.598416
$ws.Range("I4").Value = 0.3032748457625707
$ws.Range("J4").Value = 0.3032748457625707
$ws.Range("M4").Value = 0.1546876666666667
$ws.Range("N4").Value = 0.464063
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.3917937471342223
$ws.Range("R4").Value = 3.526143724208
$ws.Range("S4").Value = 0.3032748457625707
$ws.Range("T4").Value = 0.3032748457625707

# --- Remove old rows 5-7 (MuSCs->ECs, MuSCs->FAPs duplicates from before) ---
$ws.Rows("5:7").Delete()
